{"js": "// Replace the 25 \"NNN\u00d7N=\" practice-problem expressions in the single\n// table with their updated values. Each tuple is (row, col, oldText,\n// newText) using 0-based Table.getCell(row, col) addressing, listed in\n// document order (row-major across the 5 populated rows / 5 columns).\nconst replacements = [\n  [0, 0, \"204\u00d76=\", \"307\u00d73=\"],\n  [0, 1, \"866\u00d79=\", \"279\u00d78=\"],\n  [0, 2, \"298\u00d78=\", \"915\u00d76=\"],\n  [0, 3, \"333\u00d78=\", \"267\u00d77=\"],\n  [0, 4, \"461\u00d76=\", \"364\u00d78=\"],\n  [4, 0, \"905\u00d79=\", \"943\u00d79=\"],\n  [4, 1, \"110\u00d75=\", \"757\u00d73=\"],\n  [4, 2, \"636\u00d78=\", \"651\u00d73=\"],\n  [4, 3, \"430\u00d78=\", \"803\u00d72=\"],\n  [4, 4, \"505\u00d76=\", \"176\u00d77=\"],\n  [9, 0, \"556\u00d78=\", \"855\u00d77=\"],\n  [9, 1, \"141\u00d77=\", \"648\u00d78=\"],\n  [9, 2, \"623\u00d75=\", \"795\u00d77=\"],\n  [9, 3, \"309\u00d76=\", \"275\u00d74=\"],\n  [9, 4, \"383\u00d77=\", \"840\u00d72=\"],\n  [14, 0, \"768\u00d73=\", \"525\u00d72=\"],\n  [14, 1, \"312\u00d76=\", \"436\u00d78=\"],\n  [14, 2, \"283\u00d75=\", \"990\u00d79=\"],\n  [14, 3, \"470\u00d76=\", \"822\u00d76=\"],\n  [14, 4, \"364\u00d78=\", \"611\u00d72=\"],\n  [19, 0, \"390\u00d72=\", \"222\u00d75=\"],\n  [19, 1, \"921\u00d72=\", \"796\u00d76=\"],\n  [19, 2, \"641\u00d75=\", \"804\u00d75=\"],\n  [19, 3, \"435\u00d72=\", \"803\u00d76=\"],\n  [19, 4, \"539\u00d76=\", \"836\u00d74=\"],\n];\n\nconst table = context.document.body.tables.getFirst();\n\n// Queue a load for the current text of every cell we're about to touch.\nconst cells = replacements.map(([row, col]) => {\n  const cell = table.getCell(row, col);\n  cell.load(\"value\");\n  return cell;\n});\nawait context.sync();\n\n// Verify, then overwrite each cell's text in place (preserves formatting).\nfor (let i = 0; i < cells.length; i++) {\n  const [, , oldText, newText] = replacements[i];\n  if (cells[i].value !== oldText) {\n    throw new Error(\n      `Unexpected text in cell ${i} (row ${replacements[i][0]}, col ${replacements[i][1]}): ${cells[i].value}`\n    );\n  }\n  cells[i].value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"NNN\u00d7N=\" practice-problem expressions in the single\n# table with their updated values. Each tuple is (row, col, oldText,\n# newText) using 1-based Table.Cell(row, col) addressing, listed in\n# document order (row-major across the 5 populated rows / 5 columns).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @(1, 1, \"204\u00d76=\", \"307\u00d73=\"),\n    @(1, 2, \"866\u00d79=\", \"279\u00d78=\"),\n    @(1, 3, \"298\u00d78=\", \"915\u00d76=\"),\n    @(1, 4, \"333\u00d78=\", \"267\u00d77=\"),\n    @(1, 5, \"461\u00d76=\", \"364\u00d78=\"),\n    @(5, 1, \"905\u00d79=\", \"943\u00d79=\"),\n    @(5, 2, \"110\u00d75=\", \"757\u00d73=\"),\n    @(5, 3, \"636\u00d78=\", \"651\u00d73=\"),\n    @(5, 4, \"430\u00d78=\", \"803\u00d72=\"),\n    @(5, 5, \"505\u00d76=\", \"176\u00d77=\"),\n    @(10, 1, \"556\u00d78=\", \"855\u00d77=\"),\n    @(10, 2, \"141\u00d77=\", \"648\u00d78=\"),\n    @(10, 3, \"623\u00d75=\", \"795\u00d77=\"),\n    @(10, 4, \"309\u00d76=\", \"275\u00d74=\"),\n    @(10, 5, \"383\u00d77=\", \"840\u00d72=\"),\n    @(15, 1, \"768\u00d73=\", \"525\u00d72=\"),\n    @(15, 2, \"312\u00d76=\", \"436\u00d78=\"),\n    @(15, 3, \"283\u00d75=\", \"990\u00d79=\"),\n    @(15, 4, \"470\u00d76=\", \"822\u00d76=\"),\n    @(15, 5, \"364\u00d78=\", \"611\u00d72=\"),\n    @(20, 1, \"390\u00d72=\", \"222\u00d75=\"),\n    @(20, 2, \"921\u00d72=\", \"796\u00d76=\"),\n    @(20, 3, \"641\u00d75=\", \"804\u00d75=\"),\n    @(20, 4, \"435\u00d72=\", \"803\u00d76=\"),\n    @(20, 5, \"539\u00d76=\", \"836\u00d74=\")\n)\n\nforeach ($rep in $replacements) {\n    $row = $rep[0]\n    $col = $rep[1]\n    $oldText = $rep[2]\n    $newText = $rep[3]\n    $cell = $t.Cell($row, $col)\n    $current = $cell.Range.Text\n    $expected = $oldText + [char]13 + [char]7\n    if ($current -ne $expected) {\n        throw (\"Unexpected text at row \" + $row + \" col \" + $col + \": \" + $current)\n    }\n    $cell.Range.Text = $newText\n}\n"}
